# Insert a new price-observation row at row 35 (weekly update adds a new
# record for this market/product and shifts the existing history down).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 45076
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100107
$ws.Range("H35").Value = "Otros"
$ws.Range("I35").Value = 100107011
$ws.Range("J35").Value = "Tuna"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Especial"
$ws.Range("M35").Value = 50
$ws.Range("N35").Value = 25000
$ws.Range("O35").Value = 25000
$ws.Range("P35").Value = 25000
$ws.Range("Q35").Value = "$/caja 16 kilos"
$ws.Range("R35").Value = "Provincia de Los Andes"
$ws.Range("S35").Value = 1562
$ws.Range("T35").Value = 16
